$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column J (rows 4-33) with the same pattern as F:I, i.e. =IF(J$2=$Cn,1,0)
$ws.Range("J4:J33").Formula = "=IF(J`$2=`$C4,1,0)"

# Update the selection to match the new sqref range
$null = $ws.Range("D2:N5").Select()
